# #5: insurance, claim, debt, investment done
#
# The "保險" (insurance) sheet (sheet9) had a buggy header row (B1:D1 were
# pointing at stray data strings instead of real header labels) and was
# missing the standard metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index) that every
# other property-type sheet in this workbook already carries.
#
# This script:
#   1. Fixes the B1:D1 header labels on the "保險" sheet to
#      "company" / "name" / "owner".
#   2. Adds the missing E:K columns (header + 5 data rows) with the
#      standard metadata, matching the layout used on the other sheets
#      (e.g. "具有相當價值之財產").
#   3. Mirrors the one incidental relabelling on the "具有相當價值之財產"
#      sheet (F2: otherbonds -> antique) that results from the shared
#      string table being edited upstream.

$wb  = $excel.ActiveWorkbook
$ws9 = $wb.Worksheets.Item(9)   # 保險
$ws8 = $wb.Worksheets.Item(8)   # 具有相當價值之財產

# --- incidental relabel on the neighbouring sheet ---------------------
$ws8.Range("F2").Value = "antique"

# --- fix the mislabeled header cells on 保險 ---------------------------
$ws9.Range("B1").Value = "company"
$ws9.Range("C1").Value = "name"
$ws9.Range("D1").Value = "owner"

# --- add the new header cells E1:K1 ------------------------------------
$newCols = @("E", "F", "G", "H", "I", "J", "K")
$headerVals = @{
    E = "property_category"
    F = "category"
    G = "date"
    H = "legislator_name"
    I = "legislator_id"
    J = "source_file"
    K = "index"
}

foreach ($col in $newCols) {
    $cell = $ws9.Range($col + "1")
    # Force text so header labels never get reinterpreted as numbers/dates.
    $cell.NumberFormat = "@"
    $cell.Value = $headerVals[$col]
    # Re-use the existing header style (same as B1:D1) instead of leaving
    # the temporary text format applied.
    $ws9.Range("D1").Copy()
    $cell.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- add the new data cells E2:K6 --------------------------------------
$rowVals = @{
    2 = @{ E = "insurance"; F = "normal"; G = "2012-05-01"; H = "黃偉哲"; I = 1367; J = "tmp62651"; K = 136 }
    3 = @{ E = "insurance"; F = "normal"; G = "2012-05-01"; H = "黃偉哲"; I = 1367; J = "tmp62651"; K = 137 }
    4 = @{ E = "insurance"; F = "normal"; G = "2012-05-01"; H = "黃偉哲"; I = 1367; J = "tmp62651"; K = 139 }
    5 = @{ E = "insurance"; F = "normal"; G = "2012-05-01"; H = "黃偉哲"; I = 1367; J = "tmp62651"; K = 140 }
    6 = @{ E = "insurance"; F = "normal"; G = "2012-05-01"; H = "黃偉哲"; I = 1367; J = "tmp62651"; K = 141 }
}

foreach ($r in $rowVals.Keys) {
    $data = $rowVals[$r]
    foreach ($col in $newCols) {
        $cell = $ws9.Range($col + $r)
        $val = $data[$col]
        if ($val -is [string]) {
            # Force text so values like the date "2012-05-01" stay literal
            # strings instead of being auto-converted to a date serial.
            $cell.NumberFormat = "@"
        }
        $cell.Value = $val
        # Re-use the existing data-row style (same as B:D on this row).
        $ws9.Range("D" + $r).Copy()
        $cell.PasteSpecial(-4122)
    }
}
$excel.CutCopyMode = $false
